# Apply updated crypto price/volume data per commit 'Updated cryptos list'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'" + '29.838.27'
$ws.Cells.Item(2, 5).Value = '  -1.59%  '
$ws.Cells.Item(3, 4).Value = "'" + '1.893.07'
$ws.Cells.Item(3, 5).Value = '  -1.32%  '
$ws.Cells.Item(4, 4).Value = "'" + '0.9999'
$ws.Cells.Item(4, 5).Value = '  -0.14%  '
$ws.Cells.Item(5, 4).Value = "'" + '0.7784'
$ws.Cells.Item(5, 5).Value = '  -3.64%  '
$ws.Cells.Item(6, 4).Value = "'" + '244.20'
$ws.Cells.Item(6, 5).Value = '  -0.04%  '
$ws.Cells.Item(7, 4).Value = "'" + '0.9996'
$ws.Cells.Item(7, 5).Value = '  -0.18%  '
$ws.Cells.Item(8, 4).Value = "'" + '0.3129'
$ws.Cells.Item(8, 5).Value = '  -3.29%  '
$ws.Cells.Item(9, 4).Value = "'" + '25.43'
$ws.Cells.Item(9, 5).Value = '  -6.25%  '
$ws.Cells.Item(10, 4).Value = "'" + '0.07202'
$ws.Cells.Item(10, 5).Value = '  +1.44%  '
$ws.Cells.Item(11, 4).Value = "'" + '0.08079'
$ws.Cells.Item(12, 4).Value = "'" + '0.7681'
$ws.Cells.Item(12, 5).Value = '  -2.13%  '
$ws.Cells.Item(13, 4).Value = "'" + '5.501'
$ws.Cells.Item(13, 5).Value = '  +1.49%  '
$ws.Cells.Item(14, 4).Value = "'" + '1.889.23'
$ws.Cells.Item(14, 5).Value = '  -1.60%  '
$ws.Cells.Item(15, 4).Value = "'" + '92.45'
$ws.Cells.Item(15, 5).Value = '  -2.30%  '
$ws.Cells.Item(16, 4).Value = "'" + '6.180'
$ws.Cells.Item(16, 5).Value = '  +2.58%  '
$ws.Cells.Item(17, 4).Value = "'" + '29.843.28'
$ws.Cells.Item(17, 5).Value = '  -1.57%  '
$ws.Cells.Item(18, 4).Value = "'" + '13.97'
$ws.Cells.Item(18, 5).Value = '  -2.31%  '
$ws.Cells.Item(19, 4).Value = "'" + '243.83'
$ws.Cells.Item(19, 5).Value = '  -2.81%  '
$ws.Cells.Item(20, 4).Value = "'" + '0.000007776'
$ws.Cells.Item(20, 5).Value = '  -0.30%  '
$ws.Cells.Item(21, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(21, 4).Value = "'" + '2.166.75'
$ws.Cells.Item(21, 5).Value = '  +0.46%  '
$ws.Cells.Item(22, 2).Value = 'Dai'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(22, 4).Value = "'" + '1.000'
$ws.Cells.Item(22, 5).Value = '  -0.06%  '
$ws.Cells.Item(23, 2).Value = 'Chainlink'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(23, 4).Value = "'" + '8.140'
$ws.Cells.Item(23, 5).Value = '  +2.04%  '
$ws.Cells.Item(24, 4).Value = "'" + '1.000'
$ws.Cells.Item(24, 5).Value = '  -0.17%  '
$ws.Cells.Item(25, 4).Value = "'" + '0.1554'
$ws.Cells.Item(25, 5).Value = '  -4.16%  '
$ws.Cells.Item(26, 4).Value = "'" + '9.401'
$ws.Cells.Item(26, 5).Value = '  -1.09%  '
$ws.Cells.Item(27, 4).Value = "'" + '162.52'
$ws.Cells.Item(27, 5).Value = '  -3.02%  '
$ws.Cells.Item(28, 4).Value = "'" + '18.74'
$ws.Cells.Item(28, 5).Value = '  -1.84%  '
$ws.Cells.Item(29, 4).Value = "'" + '2.048'
$ws.Cells.Item(29, 5).Value = '  -4.13%  '
$ws.Cells.Item(30, 4).Value = "'" + '1.430'
$ws.Cells.Item(30, 5).Value = '  +3.91%  '
$ws.Cells.Item(31, 4).Value = "'" + '1.549'
$ws.Cells.Item(31, 5).Value = '  +0.72%  '
$ws.Cells.Item(32, 4).Value = "'" + '4.474'
$ws.Cells.Item(32, 5).Value = '  +2.84%  '
$ws.Cells.Item(33, 4).Value = "'" + '4.106'
$ws.Cells.Item(33, 5).Value = '  -0.76%  '
$ws.Cells.Item(34, 4).Value = "'" + '0.05528'
$ws.Cells.Item(34, 5).Value = '  -1.26%  '
$ws.Cells.Item(35, 5).Value = '  -2.87%  '
$ws.Cells.Item(36, 4).Value = "'" + '0.7485'
$ws.Cells.Item(36, 5).Value = '  +0.60%  '
$ws.Cells.Item(37, 4).Value = "'" + '1.004'
$ws.Cells.Item(37, 5).Value = '  +0.38%  '
$ws.Cells.Item(38, 5).Value = '  -3.27%  '
$ws.Cells.Item(39, 4).Value = "'" + '0.01919'
$ws.Cells.Item(39, 5).Value = '  -1.65%  '
$ws.Cells.Item(40, 5).Value = '  -1.46%  '
$ws.Cells.Item(41, 4).Value = "'" + '1.136.40'
$ws.Cells.Item(41, 5).Value = '  +9.61%  '
$ws.Cells.Item(42, 4).Value = "'" + '73.63'
$ws.Cells.Item(42, 5).Value = '  -0.13%  '
$ws.Cells.Item(43, 4).Value = "'" + '0.4423'
$ws.Cells.Item(43, 5).Value = '  -1.29%  '
$ws.Cells.Item(44, 4).Value = "'" + '5.895'
$ws.Cells.Item(44, 5).Value = '  -1.43%  '
$ws.Cells.Item(45, 4).Value = "'" + '0.8508'
$ws.Cells.Item(45, 5).Value = '  -0.62%  '
$ws.Cells.Item(46, 2).Value = 'PaxDollar'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(46, 4).Value = "'" + '0.9995'
$ws.Cells.Item(46, 5).Value = '  -0.17%  '
$ws.Cells.Item(47, 2).Value = 'Quant'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(47, 4).Value = "'" + '103.84'
$ws.Cells.Item(47, 5).Value = '  +0.81%  '
$ws.Cells.Item(48, 4).Value = "'" + '1.892'
$ws.Cells.Item(48, 5).Value = '  -2.13%  '
$ws.Cells.Item(49, 4).Value = "'" + '9.922'
$ws.Cells.Item(49, 5).Value = '  -0.45%  '
$ws.Cells.Item(50, 4).Value = "'" + '3.042'
$ws.Cells.Item(50, 5).Value = '  +11.57%  '
$ws.Cells.Item(51, 4).Value = "'" + '7.473'
$ws.Cells.Item(51, 5).Value = '  -2.29%  '
